# Restructure the "Storage" worksheet from a vertical parameter/value list
# into a horizontal table (one row per storage technology), matching the
# layout used on the "Technologies" sheet, and add a "Pit" storage entry.

$wb = $excel.ActiveWorkbook
$storage = $wb.Worksheets.Item("Storage")

# Capture the old parameter/value pairs (rows 1-11, columns A/B) before
# clearing the sheet, so we know which value goes under which header.
$oldHeaderParam = $storage.Cells.Item(1, 1).Value2   # "parameter"
$oldHeaderValue = $storage.Cells.Item(1, 2).Value2   # "value"

$capacityCost      = $storage.Cells.Item(2, 2).Value2   # capacity_cost   25000
$fixedOm           = $storage.Cells.Item(3, 2).Value2   # fixed_om        500
$variableOm        = $storage.Cells.Item(4, 2).Value2   # variable_om     0.1
$efficiency        = $storage.Cells.Item(5, 2).Value2   # efficiency      0.95
$lossRate          = $storage.Cells.Item(6, 2).Value2   # loss_rate       0.02
$maxChargeRate     = $storage.Cells.Item(7, 2).Value2   # max_charge_rate 0.25
$maxDischargeRate  = $storage.Cells.Item(8, 2).Value2   # max_discharge_rate 0.25
$lifetime          = $storage.Cells.Item(9, 2).Value2   # lifetime        4
$maxCapacity       = $storage.Cells.Item(10, 2).Value2  # max_capacity    1000
$initialCapacity   = $storage.Cells.Item(11, 2).Value2  # initial_capacity 0

# Clear the old content entirely so the sheet ends up exactly A1:K2.
$usedRange = $storage.UsedRange
$usedRange.Clear() | Out-Null

# New header row.
$storage.Cells.Item(1, 1).Value = "storage"
$storage.Cells.Item(1, 2).Value = "capacity_cost"
$storage.Cells.Item(1, 3).Value = "fixed_om"
$storage.Cells.Item(1, 4).Value = "variable_om"
$storage.Cells.Item(1, 5).Value = "efficiency"
$storage.Cells.Item(1, 6).Value = "loss_rate"
$storage.Cells.Item(1, 7).Value = "max_charge_rate"
$storage.Cells.Item(1, 8).Value = "max_discharge_rate"
$storage.Cells.Item(1, 9).Value = "lifetime"
$storage.Cells.Item(1, 10).Value = "max_capacity"
$storage.Cells.Item(1, 11).Value = "initial_capacity"

# New data row for the "Pit" thermal storage technology.
$storage.Cells.Item(2, 1).Value = "Pit"
$storage.Cells.Item(2, 2).Value = $capacityCost
$storage.Cells.Item(2, 3).Value = $fixedOm
$storage.Cells.Item(2, 4).Value = $variableOm
$storage.Cells.Item(2, 5).Value = $efficiency
$storage.Cells.Item(2, 6).Value = $lossRate
$storage.Cells.Item(2, 7).Value = $maxChargeRate
$storage.Cells.Item(2, 8).Value = $maxDischargeRate
$storage.Cells.Item(2, 9).Value = $lifetime
$storage.Cells.Item(2, 10).Value = $maxCapacity
$storage.Cells.Item(2, 11).Value = $initialCapacity

# Auto-fit the new columns, similar to the bestFit widths seen on other
# sheets after the edit.
$storage.Columns.Item("B:K").AutoFit() | Out-Null

# Make "Storage" the active / selected sheet, and select a cell out past
# the data (mirroring the author's click on J35) as the last interaction.
$storage.Activate()
$storage.Range("J35").Select() | Out-Null
